$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data contained two duplicate rows (original rows 2 and 5) for
# the 케이엔알시스템 / DB+NH underwriting entry. Remove both entire rows,
# shifting the remaining rows up. Delete from the bottom first so the
# row numbers of rows not yet deleted remain stable.
$ws.Rows.Item(5).EntireRow.Delete() | Out-Null
$ws.Rows.Item(2).EntireRow.Delete() | Out-Null
